$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update task description text in E3 (shared string index 13).
# Leading apostrophe preserves the cell's existing quote-prefix ("text") style
# instead of Excel normalizing it away when the value is (re)written.
$ws.Range("E3").Value = "'Reception et analyse du cahier des charges, mise en place des documents de base, recherche d'informations par rapport à la puissance du moteur et des LEDs"

# Fill in new row 4 data
$ws.Range("D3").Value = 5
$ws.Range("D4").Value = 8
$ws.Range("E4").Value = "Correction du planning, recherche d'informations et de composants pour le driver de moteur pas à pas, meeting d'une heure"
$ws.Range("F4").Value = 3

# Update selection to E19
$ws.Range("E19").Select()
